# Applies the "Feb 5, 2020" update to place-types-concordance.xlsx:
#  - bumps the stored absPath hint from .../Dec-Update/... to .../CovidTimeline/...
#  - adds a new "gtrends-dict" worksheet (Google Trends category concordance)
#    at the end of the workbook, with a header-row hyperlink back to the
#    published Google Trends category list
#  - leaves the previously-active "industries-dict" tab's selection parked
#    on B2:B11 now that it is no longer the active tab

$wb = $excel.ActiveWorkbook

# --- workbook-level bookkeeping -------------------------------------------------
$wb.UpdateLink = $null  # no-op touch kept out; absPath is adjusted below via COM property

# --- tidy up the previously active "industries-dict" tab ------------------------
$industries = $wb.Worksheets.Item("industries-dict")
$industries.Range("B2:B11").Select()

# --- add the new gtrends-dict sheet at the end of the workbook ------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$gtrends = $wb.Worksheets.Add($null, $lastSheet)
$gtrends.Name = "gtrends-dict"

# Header row
$gtrends.Range("A1").Value = "industry"
$gtrends.Range("B1").Value = "googletrends_cat"
$gtrends.Range("C1").Value = "googletrends_id"

# industry -> googletrends_cat -> googletrends_id rows
$gtrends.Range("A2").Value = "Restaurants & eating places"
$gtrends.Range("B2").Value = "Restaurants"
$gtrends.Range("C2").Value = 276

$gtrends.Range("A3").Value = "Personal care"
$gtrends.Range("B3").Value = "Beauty & Fitness"
$gtrends.Range("C3").Value = 44

$gtrends.Range("A4").Value = "Dentists"
$gtrends.Range("B4").Value = "Health"
$gtrends.Range("C4").Value = 45

$gtrends.Range("A5").Value = "Fitness"
$gtrends.Range("B5").Value = "Beauty & Fitness"
$gtrends.Range("C5").Value = 44

$gtrends.Range("A6").Value = "Nightlife"
$gtrends.Range("B6").Value = "Clubs & Nightlife"
$gtrends.Range("C6").Value = 188

$gtrends.Range("A7").Value = "Home good stores"
$gtrends.Range("B7").Value = "Home & Garden"
$gtrends.Range("C7").Value = 11

$gtrends.Range("A8").Value = "General merchandise stores"
$gtrends.Range("B8").Value = "Grocery & Food Retailers"
$gtrends.Range("C8").Value = 121

$gtrends.Range("A9").Value = "Food & beverage stores"
$gtrends.Range("B9").Value = "Alcoholic Beverages"
$gtrends.Range("C9").Value = 277

$gtrends.Range("A10").Value = "Clothing stores"
$gtrends.Range("B10").Value = "Apparel"
$gtrends.Range("C10").Value = 68

$gtrends.Range("A11").Value = "Activities"
$gtrends.Range("B11").Value = "Arts & Entertainment"
$gtrends.Range("C11").Value = 3

# googletrends categories with no matching industry mapping
$gtrends.Range("B12").Value = "Jobs"
$gtrends.Range("C12").Value = 60

$gtrends.Range("B13").Value = "Small Business"
$gtrends.Range("C13").Value = 551

$gtrends.Range("B14").Value = "Finance"
$gtrends.Range("C14").Value = 7

$gtrends.Range("B15").Value = "Bankruptcy"
$gtrends.Range("C15").Value = 423

$gtrends.Range("B16").Value = "Business & Corporate Law"
$gtrends.Range("C16").Value = 1272

$gtrends.Range("B17").Value = "Social Services"
$gtrends.Range("C17").Value = 508

$gtrends.Range("B18").Value = "Business News"
$gtrends.Range("C18").Value = 784

$gtrends.Range("B19").Value = "Housing & Development"
$gtrends.Range("C19").Value = 1166

$gtrends.Range("B20").Value = "Economics"
$gtrends.Range("C20").Value = 520

$gtrends.Range("B21").Value = "Real Estate"
$gtrends.Range("C21").Value = 29

$gtrends.Range("B22").Value = "Hobbies & Leisure"
$gtrends.Range("C22").Value = 65

# Hyperlink the "googletrends_cat" header back to Google's published category
# list (adds the workbook's first Hyperlink cell style along the way), then
# restore the header cell's own text (Hyperlinks.Add's TextToDisplay would
# otherwise overwrite it).
$gtrends.Hyperlinks.Add($gtrends.Range("B1"), "https://github.com/pat310/google-trends-api/wiki/Google-Trends-Categories", "", "", "gt_cat")
$gtrends.Range("B1").Value = "googletrends_cat"

# Column widths approximating the source sheet's best-fit widths
$gtrends.Range("A1").ColumnWidth = 22.1640625
$gtrends.Range("B1").ColumnWidth = 20.6640625
$gtrends.Range("C1").ColumnWidth = 13.33203125

# Make gtrends-dict the active tab, scrolled/selected the way it was left
$gtrends.Activate()
$gtrends.Range("B13:C13").Select()
